$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column F (dSF) by row number.
# Rows 15 and 19 are unchanged (remain 0).
$updates = @{
    2  = -1
    3  = 4
    4  = -4
    5  = -1
    6  = 3
    7  = -1
    8  = 12
    9  = 1
    10 = -1
    11 = 2
    12 = 0
    13 = -2
    14 = -1
    16 = 1
    17 = -2
    18 = -2
    20 = -1
    21 = 0
    22 = -2
    23 = -4
    24 = 6
    25 = 5
    26 = -2
    27 = -5
    28 = 3
    29 = 8
    30 = 3
    31 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
